# "fixed inner cable remove"
#
# Appends 8 new pull-sheet rows (sheet rows 24-31) to the bottom of the
# "Sheet1" table, continuing the existing Pull #/Local-Express/From/To/
# SK-#/Cable Type/Cable Size/length columns (A:I). The new rows reuse the
# "LOCAL"/"543+00"/"554+90"/"7C#14" combo already used elsewhere in the
# sheet for row 24, then introduce a new "EXPRESS" / "500+00"-"600+00" /
# "STAR QUAD" / "2C#6" run for rows 25-31 (closing out with the existing
# "3C#6" cable size on the last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns: A Pull#, B Local/Express, C From, D To, E SK-#, F Cable Type,
#          G Cable Size, H Len. Before Mess., I Len. After Mess.
$newRows = @(
    @(23, "LOCAL",   "543+00", "554+90", 1, "7C#14",     "PK", 0, 0),
    @(24, "EXPRESS", "500+00", "600+00", 1, "STAR QUAD", "PK", 0, 0),
    @(25, "EXPRESS", "500+00", "600+00", 2, "STAR QUAD", "PK", 0, 0),
    @(26, "EXPRESS", "500+00", "600+00", 2, "STAR QUAD", "PK", 0, 0),
    @(27, "EXPRESS", "500+00", "600+00", 2, "2C#6",      "PK", 0, 0),
    @(28, "EXPRESS", "500+00", "600+00", 2, "2C#6",      "PK", 0, 0),
    @(29, "EXPRESS", "500+00", "600+00", 2, "2C#6",      "PK", 0, 0),
    @(30, "EXPRESS", "500+00", "600+00", 2, "3C#6",      "PK", 0, 0)
)

$r = 24
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $r++
}

# The new "STAR QUAD" label is wider than anything previously in column F,
# so it no longer fits the default width -- give it an explicit best-fit
# width (matches what Excel computes when auto-fitting this column).
$ws.Columns.Item(6).ColumnWidth = 11.140625

# Scroll the window down so row 7 is at the top and leave the new last
# row selected, matching where editing finished.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("L29").Select()
